$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status text "in progress" -> "skipped" for every cell that uses it
$ws.Range("G22").Value = "skipped"
$ws.Range("G23").Value = "skipped"

# Update the active selection to reflect G24 as the last selected cell
$ws.Range("G24").Select() | Out-Null
